$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 1520
$ws.Range("J7").Value = 431.66666
$ws.Range("L7").Value = 431.66666
$ws.Range("N7").Value = -655.66666
$ws.Range("H11").Value = 19
$ws.Range("I11").Value = 19
$ws.Range("K11").Value = 19
$ws.Range("M11").Value = 121
$ws.Range("H14").Value = 1520
$ws.Range("J14").Value = 431.66666
$ws.Range("L14").Value = 431.66666
$ws.Range("N14").Value = -813.66666
$ws.Range("H101").Value = 361
$ws.Range("I101").Value = 308
$ws.Range("K101").Value = 924
$ws.Range("M101").Value = 698
$ws.Range("H132").Value = 1868.3871
$ws.Range("I132").Value = 1087.5
$ws.Range("K132").Value = 3262.5
$ws.Range("M132").Value = -732.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 18000974
$ws.Range("I6").Value = 17144248
$ws.Range("K6").Value = 17144248
$ws.Range("M6").Value = -17144075
$ws.Range("H11").Value = 7500
$ws.Range("J11").Value = 7500
$ws.Range("L11").Value = 7500
$ws.Range("N11").Value = -7788
$ws.Range("H45").Value = 1953.9
$ws.Range("I45").Value = 1132.6154
$ws.Range("J45").Value = 3479.1428
$ws.Range("K45").Value = 1132.6154
$ws.Range("L45").Value = 3479.1428
$ws.Range("M45").Value = -755.6153999999999
$ws.Range("N45").Value = -4233.1428
$ws.Range("H93").Value = 25000
$ws.Range("J93").Value = 25000
$ws.Range("L93").Value = 25000
$ws.Range("N93").Value = -29992

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H16").Value = 1100
$ws.Range("I16").Value = 1200
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 1200
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -1030
$ws.Range("N16").Value = -1340
$ws.Range("H31").Value = 8000
$ws.Range("I31").Value = 8000
$ws.Range("K31").Value = 8000
$ws.Range("M31").Value = -7748
$ws.Range("H86").Value = 3061.2307
$ws.Range("I86").Value = 1570.6316
$ws.Range("J86").Value = 7107.143
$ws.Range("K86").Value = 1570.6316
$ws.Range("L86").Value = 7107.143
$ws.Range("M86").Value = -447.6315999999999
$ws.Range("N86").Value = -9353.143
$ws.Range("H89").Value = 3061.2307
$ws.Range("I89").Value = 1570.6316
$ws.Range("J89").Value = 7107.143
$ws.Range("K89").Value = 7853.157999999999
$ws.Range("L89").Value = 35535.715
$ws.Range("M89").Value = -2237.157999999999
$ws.Range("N89").Value = -46767.715
$ws.Range("H105").Value = 1373.3125
$ws.Range("I105").Value = 1259.3334
$ws.Range("K105").Value = 1259.3334
$ws.Range("M105").Value = 487.6666
$ws.Range("H134").Value = 1560.2858
$ws.Range("J134").Value = 3065.3333
$ws.Range("L134").Value = 9195.999899999999
$ws.Range("N134").Value = -14265.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 1000
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 1000
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 1000
$ws.Range("N6").Value = -1226
$ws.Range("H7").Value = 4152.4
$ws.Range("I7").Value = 5985.1763
$ws.Range("J7").Value = 257.75
$ws.Range("K7").Value = 5985.1763
$ws.Range("L7").Value = 257.75
$ws.Range("M7").Value = -5872.1763
$ws.Range("N7").Value = -483.75
$ws.Range("H17").Value = 12999
$ws.Range("J17").Value = 12999
$ws.Range("L17").Value = 12999
$ws.Range("N17").Value = -13347
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("H28").Value = 12397
$ws.Range("J28").Value = 12397
$ws.Range("L28").Value = 12397
$ws.Range("N28").Value = -12887
$ws.Range("H31").Value = 5226.1953
$ws.Range("I31").Value = 3085.2144
$ws.Range("K31").Value = 3085.2144
$ws.Range("M31").Value = -2790.2144
$ws.Range("H34").Value = 5226.1953
$ws.Range("I34").Value = 3085.2144
$ws.Range("K34").Value = 3085.2144
$ws.Range("M34").Value = -2883.2144
$ws.Range("H41").Value = 38721.75
$ws.Range("I41").Value = 8254.5
$ws.Range("J41").Value = 48877.5
$ws.Range("K41").Value = 8254.5
$ws.Range("L41").Value = 48877.5
$ws.Range("M41").Value = -7826.5
$ws.Range("N41").Value = -49733.5
$ws.Range("H51").Value = 40782.8
$ws.Range("I51").Value = 6060
$ws.Range("J51").Value = 92867
$ws.Range("K51").Value = 6060
$ws.Range("L51").Value = 92867
$ws.Range("M51").Value = -5324
$ws.Range("N51").Value = -94339
$ws.Range("H59").Value = 79441
$ws.Range("J59").Value = 89254.664
$ws.Range("L59").Value = 89254.664
$ws.Range("N59").Value = -91544.664
$ws.Range("H60").Value = 61055
$ws.Range("I60").Value = 30000
$ws.Range("J60").Value = 67266
$ws.Range("K60").Value = 30000
$ws.Range("L60").Value = 67266
$ws.Range("M60").Value = -29489
$ws.Range("N60").Value = -68288
$ws.Range("H61").Value = 40782.8
$ws.Range("I61").Value = 6060
$ws.Range("J61").Value = 92867
$ws.Range("K61").Value = 6060
$ws.Range("L61").Value = 92867
$ws.Range("M61").Value = -5712
$ws.Range("N61").Value = -93563
$ws.Range("H97").Value = 55000
$ws.Range("J97").Value = 55000
$ws.Range("L97").Value = 55000
$ws.Range("N97").Value = -56982
$ws.Range("H134").Value = 1707.8572
$ws.Range("I134").Value = 1773.5
$ws.Range("J134").Value = 1314
$ws.Range("K134").Value = 5320.5
$ws.Range("L134").Value = 3942
$ws.Range("M134").Value = -2785.5
$ws.Range("N134").Value = -9012
$ws.Range("M6").ClearContents()
$ws.Range("M25").ClearContents()
$ws.Range("N25").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1106.3334
$ws.Range("I5").Value = 1133.7273
$ws.Range("K5").Value = 3401.1819
$ws.Range("M5").Value = -3289.1819
$ws.Range("H57").Value = 1621.4445
$ws.Range("I57").Value = 1299
$ws.Range("K57").Value = 3897
$ws.Range("M57").Value = -3338
$ws.Range("H132").Value = 3666.5
$ws.Range("I132").Value = 3062.25
$ws.Range("J132").Value = 4875
$ws.Range("K132").Value = 27560.25
$ws.Range("L132").Value = 43875
$ws.Range("M132").Value = -25030.25
$ws.Range("N132").Value = -48935
$ws.Range("H135").Value = 1106.3334
$ws.Range("I135").Value = 1133.7273
$ws.Range("K135").Value = 10203.5457
$ws.Range("M135").Value = -7668.545700000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 1825
$ws.Range("J22").Value = 3250
$ws.Range("L22").Value = 3250
$ws.Range("N22").Value = -4308
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("H99").Value = 8000
$ws.Range("I99").Value = 8000
$ws.Range("K99").Value = 8000
$ws.Range("M99").Value = -5754
$ws.Range("N25").ClearContents()
$ws.Range("N70").ClearContents()
$ws.Range("N73").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H116").Value = 174833.33
$ws.Range("J116").Value = 174833.33
$ws.Range("L116").Value = 174833.33
$ws.Range("N116").Value = -184011.33
$ws.Range("H122").Value = 2514.8572
$ws.Range("I122").Value = 2514.8572
$ws.Range("K122").Value = 7544.571599999999
$ws.Range("M122").Value = -5094.571599999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 6101.1763
$ws.Range("I4").Value = 8452.083000000001
$ws.Range("J4").Value = 459
$ws.Range("K4").Value = 8452.083000000001
$ws.Range("L4").Value = 459
$ws.Range("M4").Value = -8339.083000000001
$ws.Range("N4").Value = -685
$ws.Range("H6").Value = 5973.6665
$ws.Range("I6").Value = 10085
$ws.Range("J6").Value = 1862.3334
$ws.Range("K6").Value = 10085
$ws.Range("L6").Value = 1862.3334
$ws.Range("M6").Value = -9970
$ws.Range("N6").Value = -2092.3334
$ws.Range("H62").Value = 9556.333000000001
$ws.Range("I62").Value = 7001
$ws.Range("K62").Value = 7001
$ws.Range("M62").Value = -6377
$ws.Range("H65").Value = 9556.333000000001
$ws.Range("I65").Value = 7001
$ws.Range("K65").Value = 35005
$ws.Range("M65").Value = -31885
$ws.Range("H93").Value = 33389
$ws.Range("J93").Value = 33389
$ws.Range("L93").Value = 33389
$ws.Range("N93").Value = -38381
$ws.Range("H100").Value = 972.05884
$ws.Range("I100").Value = 787.5
$ws.Range("K100").Value = 1575
$ws.Range("M100").Value = -1034
$ws.Range("H123").Value = 20000
$ws.Range("J123").Value = 20000
$ws.Range("L123").Value = 20000
$ws.Range("N123").Value = -29800
